$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.287.99'
$ws.Range("E2").Value = '  +1.89%  '

$ws.Range("D3").Value = '3.626.08'
$ws.Range("E3").Value = '  +0.63%  '

$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '197.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '579.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.47%  '

$ws.Range("D7").Value = '3.622.13'
$ws.Range("E7").Value = '  +0.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.621'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.29%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("E10").Value = '  +0.63%  '

$ws.Range("E11").Value = '  +6.78%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000300'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +19.93%  '

$ws.Range("E14").Value = '  +1.90%  '

$ws.Range("D15").Value = '4.205.30'
$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("D16").Value = '3.621.87'
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("E17").Value = '  +0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.55'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.85%  '

$ws.Range("D19").Value = '68.205.15'
$ws.Range("E19").Value = '  +2.16%  '

$ws.Range("E20").Value = '  +0.91%  '

$ws.Range("E21").Value = '  +0.94%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '403.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.36%  '

$ws.Range("E23").Value = '  -2.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +19.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.25%  '

$ws.Range("E26").Value = '  +3.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +8.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +19.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.72'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '700.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +20.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.06%  '

$ws.Range("E35").Value = '  +5.42%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '64.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.64'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.416'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.21%  '

$ws.Range("E39").Value = '  +11.45%  '

$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("E41").Value = '  +20.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.13'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.24%  '

$ws.Range("E43").Value = '  +1.91%  '

$ws.Range("D44").Value = '3.143.57'
$ws.Range("E44").Value = '  +17.29%  '

$ws.Range("E45").Value = '  -0.23%  '

$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +23.11%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0424'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.88%  '

$ws.Range("E48").Value = '  +0.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.39%  '

$ws.Range("E50").Value = '  -0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.28%  '
